# Apply "Update with Correct Forecast output" change to the
# "Forecast Comparison" sheet:
#   - insert a new column B "Week_Start_Date"
#   - shorten the Week labels from "W01".."W16" to "W1".."W16"
#   - fill in the new Week_Start_Date column (stored as literal text, not a date serial)
#   - store is_holiday_week (now column J) as a real boolean instead of a number

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column before column B; this shifts ASIN, MyForecast, etc.
# one column to the right and widens the used range from I to J.
$ws.Columns.Item(2).Insert()

# New header for the inserted column.
$ws.Cells.Item(1, 2).Value = "Week_Start_Date"

# Per-row week label and week-start-date text values.
$weekStartDates = @(
    "2025-01-05",
    "2025-01-12",
    "2025-01-19",
    "2025-01-26",
    "2025-02-02",
    "2025-02-09",
    "2025-02-16",
    "2025-02-23",
    "2025-03-02",
    "2025-03-09",
    "2025-03-16",
    "2025-03-23",
    "2025-03-30",
    "2025-04-06",
    "2025-04-13",
    "2025-04-20"
)

for ($i = 0; $i -lt $weekStartDates.Length; $i++) {
    $row = $i + 2

    # "W01" -> "W1", "W02" -> "W2", ... "W16" stays "W16"
    $ws.Cells.Item($row, 1).Value = "W" + ($i + 1)

    # Write the date as literal text (leading quote forces text storage
    # instead of Excel auto-converting it to a date serial number).
    $ws.Cells.Item($row, 2).Value = "'" + $weekStartDates[$i]

    # is_holiday_week now lives in column J and should be a boolean.
    $holidayCell = $ws.Cells.Item($row, 10)
    if ($holidayCell.Value -eq 1) {
        $holidayCell.Value = $true
    } else {
        $holidayCell.Value = $false
    }
}
